$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.94764579283308, 50.08387478571661]"
$ws.Range("T2").Value = "[49.987506072905525, 50.083081268898276]"
$ws.Range("L3").Value = "[49.94607144201058, 50.11250991665803]"
$ws.Range("T3").Value = "[50.01161020344684, 50.103140881260565]"
